# Estadisticos Segundo Parcial 23 Mayo
# Update the "Estadisticos 2P" sheet with the second-partial stats,
# and refresh the "Estadisticos Final" sheet's Promedio accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": Blancos/Reprobados/Aprobados/Por_Apro/Promedio ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 1
$ws2P.Range("F2").Value = 48
$ws2P.Range("G2").Value = 97.95999999999999
$ws2P.Range("H2").Value = 8.199999999999999

# --- Sheet "Estadisticos Final": Promedio updated ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("H2").Value = 8.4
